$d = $word.ActiveDocument

# --- Change 2 FIRST: Insert continuous section break after the table ---
$r = $d.Range(1109, 1109)
$r.InsertBreak(3)
$newPara = $d.Paragraphs.Item(33)
$newPara.Alignment = 3
$sec1 = $d.Sections.Item(1)
$sec1.PageSetup.PaperSize = 12
$sec1.PageSetup.PageWidth = 595.3
$sec1.PageSetup.PageHeight = 841.9
$sec1.PageSetup.TopMargin = 72
$sec1.PageSetup.BottomMargin = 72
$sec1.PageSetup.LeftMargin = 144
$sec1.PageSetup.RightMargin = 144

# --- Change 3: Final sectPr ---
$sec2 = $d.Sections.Item($d.Sections.Count)
$sec2.PageSetup.SectionStart = 2
$sec2.PageSetup.PaperSize = 12
$sec2.PageSetup.PageWidth = 595.3
$sec2.PageSetup.PageHeight = 841.9
$sec2.PageSetup.TopMargin = 72
$sec2.PageSetup.BottomMargin = 72
$sec2.PageSetup.LeftMargin = 144
$sec2.PageSetup.RightMargin = 144

# --- Change 1 LAST: Table column widths ---
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 168.2
$t.Columns.Item(2).Width = 139.1
